# Imaging assay template: add example values to row 2, drop extra example
# rows 3-4, and rename the "Output" column from "Raw Data File" to "Data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3ASY05_Imaging")

# --- Header row: rename Output column -------------------------------------
$ws.Range("AR1").Value = "Output [Data]"

# --- Example data row (row 2): update/fill example values -----------------
$ws.Range("B2").Value  = "organelle"
$ws.Range("C2").Value  = "GO"
$ws.Range("D2").Value  = "https://bioregistry.io/GO:0043226"
$ws.Range("E2").Value  = "OperaLX spinning-disk confocal microscope (Perkin Elmer)"
$ws.Range("R2").Value  = "DPBO"
$ws.Range("S2").Value  = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_1000191"
$ws.Range("T2").Value  = "MetaXpress Software"
$ws.Range("Z2").Value  = "V6.1"

# --- Remove the now-unused extra example rows 3 and 4 ----------------------
$ws.Range("A3:A4").EntireRow.Delete()
